$d = $word.ActiveDocument

# Mapping of old -> new text values (date line + each multiplication cell)
$replacements = @(
    @("2025-08-05 Tuesday", "2025-08-06 Wednesday"),
    @("81×67=", "16×88="),
    @("17×64=", "24×19="),
    @("22×11=", "15×29="),
    @("83×47=", "92×12="),
    @("25×63=", "97×59="),
    @("65×30=", "93×73="),
    @("40×92=", "99×61="),
    @("42×64=", "99×69="),
    @("48×95=", "20×32="),
    @("67×61=", "16×50="),
    @("15×32=", "20×88="),
    @("20×30=", "69×42="),
    @("40×46=", "65×15="),
    @("87×71=", "26×88="),
    @("16×85=", "34×18="),
    @("91×87=", "12×17="),
    @("32×49=", "79×25="),
    @("69×30=", "67×36="),
    @("42×84=", "97×99="),
    @("84×20=", "47×61="),
    @("51×63=", "77×11="),
    @("29×13=", "39×92="),
    @("75×63=", "72×49="),
    @("49×85=", "72×54="),
    @("34×86=", "56×33=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
